$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a serial date value of 45205 (2023-10-06)
# for every data row (rows 2 through 116). Update it to 45206 (2023-10-07).
for ($r = 2; $r -le 116; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45205) {
        $cell.Value2 = 45206
    }
}
